$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 23 new survey response rows (46-68) ---
# Values are written column-by-column in the same left-to-right order the rows
# were originally authored in, so any new shared strings get appended in the same
# sequence as in the target workbook.

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 45184.7866319444
$ws.Range("B46").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C46").Value = 45184.7874421296
$ws.Range("C46").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D46").Value = "anonymous"
$ws.Range("F46").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G46").Value = "#17 Garlic Roast Beef;#2 The Outlaw;#11 Subway Club;"
$ws.Range("H46").Value = "15 - 20"
$ws.Range("I46").Value = "Male"

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 45185.166828703703
$ws.Range("B47").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C47").Value = 45185.167905092603
$ws.Range("C47").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D47").Value = "anonymous"
$ws.Range("F47").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G47").Value = "#99 Grand Slam Ham;#15 Titan Turkey;#20 Elite Chicken & Bacon Ranch;"
$ws.Range("H47").Value = "25 - 32"
$ws.Range("I47").Value = "Female"

$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 45185.340833333299
$ws.Range("B48").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C48").Value = 45185.341238425899
$ws.Range("C48").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D48").Value = "anonymous"
$ws.Range("F48").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G48").Value = "#1 The Philly;#6 The Boss;#19 Pickleball Club;"
$ws.Range("H48").Value = "38 - 43"
$ws.Range("I48").Value = "Male"

$ws.Range("A49").Value = 48
$ws.Range("B49").Value = 45185.598946759303
$ws.Range("B49").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C49").Value = 45185.599247685197
$ws.Range("C49").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D49").Value = "anonymous"
$ws.Range("F49").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G49").Value = "#33 Teriyaki Blitz;#16 All-Pro Sweet Onion Teriyaki;"
$ws.Range("H49").Value = "15 - 20"
$ws.Range("I49").Value = "Female"

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 45186.055277777799
$ws.Range("B50").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C50").Value = 45186.056064814802
$ws.Range("C50").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D50").Value = "anonymous"
$ws.Range("F50").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G50").Value = "#10 All-American Club;#11 Subway Club;#20 Elite Chicken & Bacon Ranch;"
$ws.Range("H50").Value = "15 - 20"
$ws.Range("I50").Value = "Male"

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 45186.325266203698
$ws.Range("B51").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C51").Value = 45186.325937499998
$ws.Range("C51").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D51").Value = "anonymous"
$ws.Range("F51").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G51").Value = "#8 The Great Garlic;"
$ws.Range("H51").Value = "15 - 20"
$ws.Range("I51").Value = "Female"

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 45187.105717592603
$ws.Range("B52").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C52").Value = 45187.106145833299
$ws.Range("C52").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D52").Value = "anonymous"
$ws.Range("F52").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G52").Value = "#15 Titan Turkey;#7 The Mexicali;#9 The Champ;"
$ws.Range("H52").Value = "25 - 32"
$ws.Range("I52").Value = "Female"

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 45188.139965277798
$ws.Range("B53").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C53").Value = 45188.140138888899
$ws.Range("C53").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D53").Value = "anonymous"
$ws.Range("F53").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G53").Value = "#15 Titan Turkey;#33 Teriyaki Blitz;"
$ws.Range("H53").Value = "25 - 32"
$ws.Range("I53").Value = "Male"

$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 45188.846168981501
$ws.Range("B54").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C54").Value = 45188.846458333297
$ws.Range("C54").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D54").Value = "anonymous"
$ws.Range("F54").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G54").Value = "#16 All-Pro Sweet Onion Teriyaki;"
$ws.Range("H54").Value = "25 - 32"
$ws.Range("I54").Value = "Male"

$ws.Range("A55").Value = 54
$ws.Range("B55").Value = 45188.998865740701
$ws.Range("B55").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C55").Value = 45188.9990972222
$ws.Range("C55").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D55").Value = "anonymous"
$ws.Range("F55").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G55").Value = "Veggie Delight;"
$ws.Range("H55").Value = "15 - 20"
$ws.Range("I55").Value = "Male"

$ws.Range("A56").Value = 55
$ws.Range("B56").Value = 45189.5770486111
$ws.Range("B56").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C56").Value = 45189.577465277798
$ws.Range("C56").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D56").Value = "anonymous"
$ws.Range("F56").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G56").Value = "#1 The Philly;#20 Elite Chicken & Bacon Ranch;#3 The Monster;"
$ws.Range("H56").Value = "15 - 20"
$ws.Range("I56").Value = "Male"

$ws.Range("A57").Value = 56
$ws.Range("B57").Value = 45191.550300925897
$ws.Range("B57").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C57").Value = 45191.550787036998
$ws.Range("C57").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D57").Value = "anonymous"
$ws.Range("F57").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G57").Value = "#23 The Hotshot Italiano;#20 Elite Chicken & Bacon Ranch;"
$ws.Range("H57").Value = "38 - 43"
$ws.Range("I57").Value = "Male"

$ws.Range("A58").Value = 57
$ws.Range("B58").Value = 45191.550821759301
$ws.Range("B58").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C58").Value = 45191.551076388903
$ws.Range("C58").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D58").Value = "anonymous"
$ws.Range("F58").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G58").Value = "#11 Subway Club;"
$ws.Range("H58").Value = "38 - 43"
$ws.Range("I58").Value = "Female"

$ws.Range("A59").Value = 58
$ws.Range("B59").Value = 45191.551111111097
$ws.Range("B59").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C59").Value = 45191.551238425898
$ws.Range("C59").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D59").Value = "anonymous"
$ws.Range("F59").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G59").Value = "Veggie Delight;"
$ws.Range("H59").Value = "38 - 43"
$ws.Range("I59").Value = "Female"

$ws.Range("A60").Value = 59
$ws.Range("B60").Value = 45191.551284722198
$ws.Range("B60").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C60").Value = 45191.551863425899
$ws.Range("C60").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D60").Value = "anonymous"
$ws.Range("F60").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G60").Value = "#1 The Philly;#33 Teriyaki Blitz;#6 The Boss;"
$ws.Range("H60").Value = "15 - 20"
$ws.Range("I60").Value = "Male"

$ws.Range("A61").Value = 60
$ws.Range("B61").Value = 45191.5518981481
$ws.Range("B61").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C61").Value = 45191.552129629599
$ws.Range("C61").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D61").Value = "anonymous"
$ws.Range("F61").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G61").Value = "#15 Titan Turkey;#33 Teriyaki Blitz;#20 Elite Chicken & Bacon Ranch;"
$ws.Range("H61").Value = "15 - 20"
$ws.Range("I61").Value = "Male"

$ws.Range("A62").Value = 61
$ws.Range("B62").Value = 45191.552152777796
$ws.Range("B62").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C62").Value = 45191.552326388897
$ws.Range("C62").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D62").Value = "anonymous"
$ws.Range("F62").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G62").Value = "#15 Titan Turkey;#6 The Boss;#20 Elite Chicken & Bacon Ranch;"
$ws.Range("H62").Value = "15 - 20"
$ws.Range("I62").Value = "Male"

$ws.Range("A63").Value = 62
$ws.Range("B63").Value = 45191.552349537
$ws.Range("B63").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C63").Value = 45191.5524421296
$ws.Range("C63").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D63").Value = "anonymous"
$ws.Range("F63").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G63").Value = "#20 Elite Chicken & Bacon Ranch;"
$ws.Range("H63").Value = "15 - 20"
$ws.Range("I63").Value = "Male"

$ws.Range("A64").Value = 63
$ws.Range("B64").Value = 45191.552465277797
$ws.Range("B64").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C64").Value = 45191.552731481497
$ws.Range("C64").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D64").Value = "anonymous"
$ws.Range("F64").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G64").Value = "#23 The Hotshot Italiano;#6 The Boss;#33 Teriyaki Blitz;"
$ws.Range("H64").Value = "15 - 20"
$ws.Range("I64").Value = "Female"

$ws.Range("A65").Value = 64
$ws.Range("B65").Value = 45191.552743055603
$ws.Range("B65").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C65").Value = 45191.552835648101
$ws.Range("C65").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D65").Value = "anonymous"
$ws.Range("F65").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G65").Value = "#33 Teriyaki Blitz;"
$ws.Range("H65").Value = "15 - 20"
$ws.Range("I65").Value = "Female"

$ws.Range("A66").Value = 65
$ws.Range("B66").Value = 45191.961064814801
$ws.Range("B66").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C66").Value = 45191.961550925902
$ws.Range("C66").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D66").Value = "anonymous"
$ws.Range("F66").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G66").Value = "#30 The Beast;#33 Teriyaki Blitz;#20 Elite Chicken & Bacon Ranch;"
$ws.Range("H66").Value = "38 - 43"
$ws.Range("I66").Value = "Male"

$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 45192.572581018503
$ws.Range("B67").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C67").Value = 45192.572812500002
$ws.Range("C67").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D67").Value = "anonymous"
$ws.Range("F67").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G67").Value = "Veggie Delight;"
$ws.Range("H67").Value = "25 - 32"
$ws.Range("I67").Value = "Male"

$ws.Range("A68").Value = 67
$ws.Range("B68").Value = 45195.453090277799
$ws.Range("B68").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C68").Value = 45195.453460648103
$ws.Range("C68").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D68").Value = "anonymous"
$ws.Range("F68").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("G68").Value = "#6 The Boss;#33 Teriyaki Blitz;"
$ws.Range("H68").Value = "25 - 32"
$ws.Range("I68").Value = "Male"

# --- Resize the Excel Table (ListObject) + AutoFilter to cover the new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:I68"))

# --- Re-unify column widths A:I to width 20 (matches bestFit columns already in use) ---
$ws.Range("A1:I1").EntireColumn.ColumnWidth = 19.17

# --- Sheet view bookkeeping: drop the frozen/scrolled topLeftCell and move the ---
# --- remembered selection from I27 to L17, matching the saved author view.    ---
$ws.Range("L17").Select()
